$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values are kept as text (matches source formatting, e.g. "27.045.48", trailing zeros)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '27.045.48'
$ws.Cells.Item(2, 5).Value = '  -2.32%  '

$ws.Cells.Item(3, 4).Value = '1.872.17'
$ws.Cells.Item(3, 5).Value = '  -1.25%  '

$ws.Cells.Item(4, 4).Value = '1.002'
$ws.Cells.Item(4, 5).Value = '  +0.15%  '

$ws.Cells.Item(5, 4).Value = '306.15'
$ws.Cells.Item(5, 5).Value = '  -1.56%  '

$ws.Cells.Item(6, 5).Value = '  +0.08%  '

$ws.Cells.Item(7, 4).Value = '0.5065'
$ws.Cells.Item(7, 5).Value = '  -2.97%  '

$ws.Cells.Item(8, 4).Value = '0.3682'
$ws.Cells.Item(8, 5).Value = '  -2.57%  '

$ws.Cells.Item(9, 4).Value = '0.07155'
$ws.Cells.Item(9, 5).Value = '  -0.96%  '

$ws.Cells.Item(10, 4).Value = '0.8892'
$ws.Cells.Item(10, 5).Value = '  -1.13%  '

$ws.Cells.Item(11, 4).Value = '20.68'
$ws.Cells.Item(11, 5).Value = '  -1.93%  '

$ws.Cells.Item(12, 4).Value = '1.896.35'
$ws.Cells.Item(12, 5).Value = '  +0.21%  '

$ws.Cells.Item(13, 4).Value = '0.07522'
$ws.Cells.Item(13, 5).Value = '  -1.42%  '

$ws.Cells.Item(14, 4).Value = '5.279'
$ws.Cells.Item(14, 5).Value = '  -2.76%  '

$ws.Cells.Item(15, 4).Value = '90.89'
$ws.Cells.Item(15, 5).Value = '  -1.00%  '

$ws.Cells.Item(16, 4).Value = '1.002'
$ws.Cells.Item(16, 5).Value = '  +0.16%  '

$ws.Cells.Item(17, 4).Value = '0.000008574'
$ws.Cells.Item(17, 5).Value = '  -1.13%  '

$ws.Cells.Item(18, 4).Value = '14.05'
$ws.Cells.Item(18, 5).Value = '  -1.34%  '

$ws.Cells.Item(19, 4).Value = '0.9999'
$ws.Cells.Item(19, 5).Value = '  -0.01%  '

$ws.Cells.Item(20, 4).Value = '27.095.32'
$ws.Cells.Item(20, 5).Value = '  -2.17%  '

$ws.Cells.Item(21, 4).Value = '5.021'
$ws.Cells.Item(21, 5).Value = '  -2.33%  '

$ws.Cells.Item(22, 4).Value = '2.119.82'
$ws.Cells.Item(22, 5).Value = '  -0.39%  '

$ws.Cells.Item(23, 4).Value = '10.39'
$ws.Cells.Item(23, 5).Value = '  -3.72%  '

$ws.Cells.Item(24, 4).Value = '6.481'
$ws.Cells.Item(24, 5).Value = '  -1.76%  '

$ws.Cells.Item(25, 4).Value = '1.834'
$ws.Cells.Item(25, 5).Value = '  -0.95%  '

$ws.Cells.Item(26, 4).Value = '146.35'
$ws.Cells.Item(26, 5).Value = '  -4.41%  '

$ws.Cells.Item(27, 4).Value = '17.91'
$ws.Cells.Item(27, 5).Value = '  -1.74%  '

$ws.Cells.Item(28, 4).Value = '2.077'
$ws.Cells.Item(28, 5).Value = '  -3.78%  '

$ws.Cells.Item(29, 4).Value = '112.99'
$ws.Cells.Item(29, 5).Value = '  -0.92%  '

$ws.Cells.Item(30, 4).Value = '4.643'
$ws.Cells.Item(30, 5).Value = '  -3.16%  '

$ws.Cells.Item(31, 4).Value = '4.706'
$ws.Cells.Item(31, 5).Value = '  -2.24%  '

$ws.Cells.Item(32, 4).Value = '0.09227'
$ws.Cells.Item(32, 5).Value = '  +1.42%  '

$ws.Cells.Item(33, 4).Value = '0.05122'
$ws.Cells.Item(33, 5).Value = '  -2.80%  '

$ws.Cells.Item(34, 4).Value = '3.087'
$ws.Cells.Item(34, 5).Value = '  -2.86%  '

$ws.Cells.Item(35, 4).Value = '1.157'
$ws.Cells.Item(35, 5).Value = '  -5.19%  '

$ws.Cells.Item(36, 4).Value = '0.7327'
$ws.Cells.Item(36, 5).Value = '  -5.11%  '

$ws.Cells.Item(37, 4).Value = '3.228'
$ws.Cells.Item(37, 5).Value = '  +5.08%  '

$ws.Cells.Item(38, 4).Value = '0.02020'
$ws.Cells.Item(38, 5).Value = '  -2.89%  '

$ws.Cells.Item(39, 4).Value = '2.476'
$ws.Cells.Item(39, 5).Value = '  -4.11%  '

$ws.Cells.Item(40, 4).Value = '1.075'
$ws.Cells.Item(40, 5).Value = '  -1.50%  '

$ws.Cells.Item(41, 4).Value = '0.5298'
$ws.Cells.Item(41, 5).Value = '  -4.54%  '

$ws.Cells.Item(42, 2).Value = 'Quant'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(42, 4).Value = '116.97'
$ws.Cells.Item(42, 5).Value = '  -0.33%  '

$ws.Cells.Item(43, 2).Value = 'FraxShare'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(43, 4).Value = '6.478'
$ws.Cells.Item(43, 5).Value = '  -2.82%  '

$ws.Cells.Item(44, 4).Value = '8.446'
$ws.Cells.Item(44, 5).Value = '  -3.08%  '

$ws.Cells.Item(45, 4).Value = '0.1475'
$ws.Cells.Item(45, 5).Value = '  -2.68%  '

$ws.Cells.Item(46, 2).Value = 'PaxDollar'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(46, 4).Value = '1.001'
$ws.Cells.Item(46, 5).Value = '  +0.06%  '

$ws.Cells.Item(47, 2).Value = 'Decentraland'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(47, 4).Value = '0.4625'
$ws.Cells.Item(47, 5).Value = '  -3.74%  '

$ws.Cells.Item(48, 4).Value = '9.990'
$ws.Cells.Item(48, 5).Value = '  -4.20%  '

$ws.Cells.Item(49, 4).Value = '1.557'
$ws.Cells.Item(49, 5).Value = '  -1.78%  '

$ws.Cells.Item(50, 4).Value = '37.08'
$ws.Cells.Item(50, 5).Value = '  +0.19%  '

$ws.Cells.Item(51, 4).Value = '63.25'
$ws.Cells.Item(51, 5).Value = '  -4.64%  '
